# Auto-generated cell updates applying the Balmung_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 368286.62
$ws.Range("J17").Value = 441330.6
$ws.Range("L17").Value = 1323991.8
$ws.Range("N17").Value = -1324327.8
$ws.Range("H40").Value = 1559.6
$ws.Range("I40").Value = 1374.75
$ws.Range("K40").Value = 1374.75
$ws.Range("M40").Value = -1199.75
$ws.Range("H98").Value = 4472.8066
$ws.Range("I98").Value = 3665.5
$ws.Range("J98").Value = 5590.615
$ws.Range("K98").Value = 3665.5
$ws.Range("L98").Value = 5590.615
$ws.Range("M98").Value = -2167.5
$ws.Range("N98").Value = -8586.615
$ws.Range("H122").Value = 4472.8066
$ws.Range("I122").Value = 3665.5
$ws.Range("J122").Value = 5590.615
$ws.Range("K122").Value = 10996.5
$ws.Range("L122").Value = 16771.845
$ws.Range("M122").Value = -8546.5
$ws.Range("N122").Value = -21671.845
$ws.Range("H132").Value = 1841
$ws.Range("I132").Value = 1712.3889
$ws.Range("K132").Value = 5137.1667
$ws.Range("M132").Value = -2607.1667
$ws.Range("H133").Value = 125000
$ws.Range("J133").Value = 125000
$ws.Range("L133").Value = 125000
$ws.Range("N133").Value = -135120
$ws.Range("H137").Value = 12501662
$ws.Range("I137").Value = 1799.5
$ws.Range("K137").Value = 5398.5
$ws.Range("M137").Value = -2848.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 53951.15
$ws.Range("J45").Value = 4770.143
$ws.Range("L45").Value = 4770.143
$ws.Range("N45").Value = -5524.143
$ws.Range("H61").Value = 843190.25
$ws.Range("I61").Value = 2656.0393
$ws.Range("K61").Value = 2656.0393
$ws.Range("M61").Value = -2444.0393
$ws.Range("H110").Value = 1480
$ws.Range("I110").Value = 973
$ws.Range("J110").Value = 2747.5
$ws.Range("K110").Value = 973
$ws.Range("L110").Value = 2747.5
$ws.Range("M110").Value = 1072
$ws.Range("N110").Value = -6837.5
$ws.Range("H122").Value = 1290.9375
$ws.Range("I122").Value = 1312.0769
$ws.Range("K122").Value = 3936.2307
$ws.Range("M122").Value = -1486.2307
$ws.Range("M132").Value = -8634.125
$ws.Range("N132").Value = $null
$ws.Range("H136").Value = 843190.25
$ws.Range("I136").Value = 2656.0393
$ws.Range("K136").Value = 7968.117899999999
$ws.Range("M136").Value = -5418.117899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 18385
$ws.Range("I105").Value = 21439.2
$ws.Range("K105").Value = 21439.2
$ws.Range("M105").Value = -19692.2
$ws.Range("H134").Value = 21953436
$ws.Range("I134").Value = 2215.2964
$ws.Range("K134").Value = 6645.889200000001
$ws.Range("M134").Value = -4110.889200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1670416.6
$ws.Range("I6").Value = 10000000
$ws.Range("J6").Value = 4500
$ws.Range("K6").Value = 10000000
$ws.Range("L6").Value = 4500
$ws.Range("M6").Value = -9999887
$ws.Range("N6").Value = -4726
$ws.Range("H7").Value = 108.63158
$ws.Range("I7").Value = 118.47059
$ws.Range("K7").Value = 118.47059
$ws.Range("M7").Value = -5.470590000000001
$ws.Range("H107").Value = 2350.5186
$ws.Range("I107").Value = 2134.2666
$ws.Range("K107").Value = 2134.2666
$ws.Range("M107").Value = -214.2665999999999
$ws.Range("H122").Value = 3555.7
$ws.Range("I122").Value = 3861.889
$ws.Range("H132").Value = 17885.555
$ws.Range("I132").Value = 21253.432
$ws.Range("J132").Value = 3572.0833
$ws.Range("K132").Value = 63760.296
$ws.Range("L132").Value = 10716.2499
$ws.Range("M132").Value = -61230.296
$ws.Range("N132").Value = -15776.2499
$ws.Range("H134").Value = 1677.4131
$ws.Range("I134").Value = 1320
$ws.Range("J134").Value = 2416.0667
$ws.Range("K134").Value = 3960
$ws.Range("L134").Value = 7248.2001
$ws.Range("M134").Value = -1425
$ws.Range("N134").Value = -12318.2001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 14850.4
$ws.Range("J63").Value = 16388
$ws.Range("L63").Value = 49164
$ws.Range("N63").Value = -50662
$ws.Range("H66").Value = 14850.4
$ws.Range("J66").Value = 16388
$ws.Range("L66").Value = 147492
$ws.Range("N66").Value = -154980
$ws.Range("H137").Value = 2614.2273
$ws.Range("I137").Value = 2115.3333
$ws.Range("J137").Value = 3683.2856
$ws.Range("K137").Value = 6345.999899999999
$ws.Range("L137").Value = 11049.8568
$ws.Range("M137").Value = -1245.999899999999
$ws.Range("N137").Value = -21249.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3723.3142
$ws.Range("I22").Value = 1229.6923
$ws.Range("J22").Value = 5196.8184
$ws.Range("K22").Value = 1229.6923
$ws.Range("L22").Value = 5196.8184
$ws.Range("M22").Value = -934.6922999999999
$ws.Range("N22").Value = -5786.8184
$ws.Range("H27").Value = 3723.3142
$ws.Range("I27").Value = 1229.6923
$ws.Range("J27").Value = 5196.8184
$ws.Range("K27").Value = 1229.6923
$ws.Range("L27").Value = 5196.8184
$ws.Range("M27").Value = -1122.6923
$ws.Range("N27").Value = -5410.8184
$ws.Range("H40").Value = 1863.8334
$ws.Range("I40").Value = 1118.8788
$ws.Range("J40").Value = 4595.3335
$ws.Range("K40").Value = 1118.8788
$ws.Range("L40").Value = 4595.3335
$ws.Range("M40").Value = -982.8788
$ws.Range("N40").Value = -4867.3335
$ws.Range("H54").Value = 31666.666
$ws.Range("I54").Value = 35000
$ws.Range("K54").Value = 35000
$ws.Range("M54").Value = -34356
$ws.Range("H61").Value = 2704.3704
$ws.Range("I61").Value = 2618.261
$ws.Range("K61").Value = 2618.261
$ws.Range("M61").Value = -2416.261
$ws.Range("H113").Value = 2704.3704
$ws.Range("I113").Value = 2618.261
$ws.Range("K113").Value = 2618.261
$ws.Range("M113").Value = -448.261
$ws.Range("H132").Value = 3326.5454
$ws.Range("I132").Value = 3120.2354
$ws.Range("K132").Value = 9360.706200000001
$ws.Range("M132").Value = -6830.706200000001
$ws.Range("H136").Value = 2363.5112
$ws.Range("I136").Value = 2572.5
$ws.Range("J136").Value = 2278.6094
$ws.Range("K136").Value = 7717.5
$ws.Range("L136").Value = 6835.8282
$ws.Range("M136").Value = -5167.5
$ws.Range("N136").Value = -11935.8282

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3379.3
$ws.Range("I62").Value = 3141
$ws.Range("K62").Value = 3141
$ws.Range("M62").Value = -2517
$ws.Range("H65").Value = 3379.3
$ws.Range("I65").Value = 3141
$ws.Range("K65").Value = 15705
$ws.Range("M65").Value = -12585
$ws.Range("H126").Value = 2551.647
$ws.Range("I126").Value = 2443.4285
$ws.Range("J126").Value = 3056.6667
$ws.Range("K126").Value = 7330.2855
$ws.Range("L126").Value = 9170.000100000001
$ws.Range("M126").Value = -4860.2855
$ws.Range("N126").Value = -14110.0001
$ws.Range("H132").Value = 1808.2941
$ws.Range("I132").Value = 1388.0769
$ws.Range("K132").Value = 4164.2307
$ws.Range("M132").Value = -1634.2307
$ws.Range("H141").Value = 97427.28999999999
$ws.Range("J141").Value = 97427.28999999999
$ws.Range("L141").Value = 97427.28999999999
$ws.Range("N141").Value = -107787.29

Write-Output "Applied 180 cell updates"
